$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal TEXT value into a cell (not auto-converted to a
# number) by building it as a formula that evaluates to a string literal,
# then converting the formula to its value in place (Copy + PasteSpecial
# values-only). This mirrors how Excel stores a typed value that happens
# to look numeric but was entered/kept as text, without leaving behind an
# extraneous quote-prefix style.
function Set-TextValue($addr, $text) {
    $c = $ws.Range($addr)
    $escaped = $text.Replace('"', '""')
    $c.Formula = '="' + $escaped + '"'
    $c.Copy() | Out-Null
    $c.PasteSpecial(-4163) | Out-Null
}

# New measurement data added in columns H (another reading), I (another
# reading) and K (phase, in degrees) for rows 2..33.
$data = @(
    @(2,  "9.88", "9.2",   -21),
    @(3,  "9.86", "7.8",   -37),
    @(4,  "9.84", "6.29",  -50),
    @(5,  "9.81", "4.98",  -60),
    @(6,  "9.80", "3.92",  -67),
    @(7,  "9.78", "3.03",  -72),
    @(8,  "9.77", "2.26",  -78),
    @(9,  "9.76", "1.6",   -82),
    @(10, "9.75", "1.05",  -86),
    @(11, "9.74", "0.57",  -90),
    @(12, "9.74", "0.34",  -95),
    @(13, "9.74", "0.26",  -96),
    @(14, "9.74", "0.2",   -99),
    @(15, "9.74", "0.166", -102),
    @(16, "9.73", "0.124", -105),
    @(17, "9.73", "0.084", -116),
    @(18, "9.73", "0.048", -136),
    @(19, "9.73", "0.038", 165),
    @(20, "9.73", "0.06",  123),
    @(21, "9.73", "0.096", 110),
    @(22, "9.73", "0.292", 95),
    @(23, "9.72", "0.65",  89),
    @(24, "9.71", "1.0",   86),
    @(25, "9.67", "2.7",   73),
    @(26, "9.59", "4.59",  60),
    @(27, "9.53", "5.79",  52),
    @(28, "9.48", "6.65",  45),
    @(29, "9.45", "7.2",   39),
    @(30, "9.42", "7.63",  34),
    @(31, "9.4",  "7.93",  30),
    @(32, "9.39", "8.16",  28),
    @(33, "9.38", "8.32",  25)
)

foreach ($row in $data) {
    $r = $row[0]
    $hVal = $row[1]
    $iVal = $row[2]
    $kVal = $row[3]

    Set-TextValue ("H" + $r) $hVal

    # Row 7's second reading carries a leftover custom (date-looking)
    # number format on the cell even though the content is plain text.
    if ($r -eq 7) {
        $ws.Range("I7").NumberFormat = "mmm-yy"
    }
    Set-TextValue ("I" + $r) $iVal

    $ws.Range("K" + $r).Value = $kVal
}

# Restore the selection to where the sheet was left (bottom of the new data).
$ws.Range("K33").Select() | Out-Null
